$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(98, 8).Value = 995.3333
$ws.Cells.Item(98, 9).Value = 1000
$ws.Cells.Item(98, 10).Value = 993
$ws.Cells.Item(98, 11).Value = 1000
$ws.Cells.Item(98, 12).Value = 993
$ws.Cells.Item(98, 13).Value = 498
$ws.Cells.Item(98, 14).Value = -3989
$ws.Cells.Item(122, 8).Value = 995.3333
$ws.Cells.Item(122, 9).Value = 1000
$ws.Cells.Item(122, 10).Value = 993
$ws.Cells.Item(122, 11).Value = 3000
$ws.Cells.Item(122, 12).Value = 2979
$ws.Cells.Item(122, 13).Value = -550
$ws.Cells.Item(122, 14).Value = -7879
$ws.Cells.Item(124, 8).Value = 43734.75
$ws.Cells.Item(124, 10).Value = 43734.75
$ws.Cells.Item(124, 12).Value = 43734.75
$ws.Cells.Item(124, 14).Value = -53554.75
$ws.Cells.Item(132, 8).Value = 1956.9395
$ws.Cells.Item(132, 9).Value = 2084.4075
$ws.Cells.Item(132, 10).Value = 1383.3334
$ws.Cells.Item(132, 11).Value = 6253.2225
$ws.Cells.Item(132, 12).Value = 4150.0002
$ws.Cells.Item(132, 13).Value = -3723.2225
$ws.Cells.Item(132, 14).Value = -9210.0002

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 16134292
$ws.Cells.Item(32, 9).Value = 26317876
$ws.Cells.Item(32, 10).Value = 10284.792
$ws.Cells.Item(32, 11).Value = 26317876
$ws.Cells.Item(32, 12).Value = 10284.792
$ws.Cells.Item(32, 13).Value = -26317589
$ws.Cells.Item(32, 14).Value = -10858.792
$ws.Cells.Item(61, 8).Value = 2619.2104
$ws.Cells.Item(61, 9).Value = 1803.7
$ws.Cells.Item(61, 10).Value = 3525.3333
$ws.Cells.Item(61, 11).Value = 1803.7
$ws.Cells.Item(61, 12).Value = 3525.3333
$ws.Cells.Item(61, 13).Value = -1591.7
$ws.Cells.Item(61, 14).Value = -3949.3333
$ws.Cells.Item(74, 8).Value = 4052.7715
$ws.Cells.Item(74, 9).Value = 5597.476
$ws.Cells.Item(74, 10).Value = 1735.7142
$ws.Cells.Item(74, 11).Value = 5597.476
$ws.Cells.Item(74, 12).Value = 1735.7142
$ws.Cells.Item(74, 13).Value = -4723.476
$ws.Cells.Item(74, 14).Value = -3483.7142
$ws.Cells.Item(77, 8).Value = 4052.7715
$ws.Cells.Item(77, 9).Value = 5597.476
$ws.Cells.Item(77, 10).Value = 1735.7142
$ws.Cells.Item(77, 11).Value = 27987.38
$ws.Cells.Item(77, 12).Value = 8678.571
$ws.Cells.Item(77, 13).Value = -23619.38
$ws.Cells.Item(77, 14).Value = -17414.571
$ws.Cells.Item(132, 8).Value = 2412.4102
$ws.Cells.Item(132, 9).Value = 2035.4839
$ws.Cells.Item(132, 11).Value = 6106.4517
$ws.Cells.Item(132, 13).Value = -3576.4517
$ws.Cells.Item(136, 8).Value = 2619.2104
$ws.Cells.Item(136, 9).Value = 1803.7
$ws.Cells.Item(136, 10).Value = 3525.3333
$ws.Cells.Item(136, 11).Value = 5411.1
$ws.Cells.Item(136, 12).Value = 10575.9999
$ws.Cells.Item(136, 13).Value = -2861.1
$ws.Cells.Item(136, 14).Value = -15675.9999

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(5, 8).Value = 940
$ws.Cells.Item(5, 9).Value = 940
$ws.Cells.Item(5, 11).Value = 940
$ws.Cells.Item(5, 13).Value = -827
$ws.Cells.Item(134, 8).Value = 3258.5
$ws.Cells.Item(134, 9).Value = 3076.5334
$ws.Cells.Item(134, 11).Value = 9229.600199999999
$ws.Cells.Item(134, 13).Value = -6694.600199999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 1487.375
$ws.Cells.Item(31, 9).Value = 1124.7142
$ws.Cells.Item(31, 10).Value = 1769.4445
$ws.Cells.Item(31, 11).Value = 1124.7142
$ws.Cells.Item(31, 12).Value = 1769.4445
$ws.Cells.Item(31, 13).Value = -829.7141999999999
$ws.Cells.Item(31, 14).Value = -2359.4445
$ws.Cells.Item(34, 8).Value = 1487.375
$ws.Cells.Item(34, 9).Value = 1124.7142
$ws.Cells.Item(34, 10).Value = 1769.4445
$ws.Cells.Item(34, 11).Value = 1124.7142
$ws.Cells.Item(34, 12).Value = 1769.4445
$ws.Cells.Item(34, 13).Value = -922.7141999999999
$ws.Cells.Item(34, 14).Value = -2173.4445
$ws.Cells.Item(38, 8).Value = 12199.75
$ws.Cells.Item(38, 9).Value = 4499.5
$ws.Cells.Item(38, 10).Value = 19900
$ws.Cells.Item(38, 11).Value = 4499.5
$ws.Cells.Item(38, 12).Value = 19900
$ws.Cells.Item(38, 13).Value = -4122.5
$ws.Cells.Item(38, 14).Value = -20654
$ws.Cells.Item(46, 8).Value = 12199.75
$ws.Cells.Item(46, 9).Value = 4499.5
$ws.Cells.Item(46, 10).Value = 19900
$ws.Cells.Item(46, 11).Value = 4499.5
$ws.Cells.Item(46, 12).Value = 19900
$ws.Cells.Item(46, 13).Value = -4288.5
$ws.Cells.Item(46, 14).Value = -20322
$ws.Cells.Item(132, 8).Value = 62502484
$ws.Cells.Item(132, 9).Value = 125001704
$ws.Cells.Item(132, 10).Value = 3265.75
$ws.Cells.Item(132, 11).Value = 375005112
$ws.Cells.Item(132, 12).Value = 9797.25
$ws.Cells.Item(132, 13).Value = -375002582
$ws.Cells.Item(132, 14).Value = -14857.25
$ws.Cells.Item(134, 8).Value = 3908.5312
$ws.Cells.Item(134, 9).Value = 2237.35
$ws.Cells.Item(134, 11).Value = 6712.049999999999
$ws.Cells.Item(134, 13).Value = -4177.049999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(131, 8).Value = 1117.7858
$ws.Cells.Item(131, 10).Value = 1143.7693
$ws.Cells.Item(131, 12).Value = 3431.3079
$ws.Cells.Item(131, 14).Value = -13511.3079

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 5962.815
$ws.Cells.Item(70, 9).Value = 6057.5386
$ws.Cells.Item(70, 11).Value = 6057.5386
$ws.Cells.Item(70, 13).Value = -5787.5386
$ws.Cells.Item(73, 8).Value = 5962.815
$ws.Cells.Item(73, 9).Value = 6057.5386
$ws.Cells.Item(73, 11).Value = 6057.5386
$ws.Cells.Item(73, 13).Value = -5121.5386
$ws.Cells.Item(122, 8).Value = 3858.9092
$ws.Cells.Item(122, 9).Value = 4449.6665
$ws.Cells.Item(122, 10).Value = 3150
$ws.Cells.Item(122, 11).Value = 13348.9995
$ws.Cells.Item(122, 12).Value = 9450
$ws.Cells.Item(122, 13).Value = -10898.9995
$ws.Cells.Item(122, 14).Value = -14350

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 3954.2856
$ws.Cells.Item(7, 9).Value = 4220
$ws.Cells.Item(7, 10).Value = 3600
$ws.Cells.Item(7, 11).Value = 4220
$ws.Cells.Item(7, 12).Value = 3600
$ws.Cells.Item(7, 13).Value = -4108
$ws.Cells.Item(7, 14).Value = -3824
$ws.Cells.Item(40, 8).Value = 5256.3
$ws.Cells.Item(40, 9).Value = 3871
$ws.Cells.Item(40, 10).Value = 6179.8335
$ws.Cells.Item(40, 11).Value = 3871
$ws.Cells.Item(40, 12).Value = 6179.8335
$ws.Cells.Item(40, 13).Value = -3735
$ws.Cells.Item(40, 14).Value = -6451.8335
$ws.Cells.Item(126, 8).Value = 3954.2856
$ws.Cells.Item(126, 9).Value = 4220
$ws.Cells.Item(126, 10).Value = 3600
$ws.Cells.Item(126, 11).Value = 12660
$ws.Cells.Item(126, 12).Value = 10800
$ws.Cells.Item(126, 13).Value = -10190
$ws.Cells.Item(126, 14).Value = -15740
$ws.Cells.Item(132, 8).Value = 12439
$ws.Cells.Item(132, 9).Value = 12555.714
$ws.Cells.Item(132, 10).Value = 12166.667
$ws.Cells.Item(132, 11).Value = 37667.142
$ws.Cells.Item(132, 12).Value = 36500.001
$ws.Cells.Item(132, 13).Value = -35137.142
$ws.Cells.Item(132, 14).Value = -41560.001
$ws.Cells.Item(133, 8).Value = 51698
$ws.Cells.Item(133, 10).Value = 51698
$ws.Cells.Item(133, 12).Value = 51698
$ws.Cells.Item(133, 14).Value = -56758

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(126, 8).Value = 1191.4706
$ws.Cells.Item(126, 9).Value = 1085.4546
$ws.Cells.Item(126, 10).Value = 1385.8334
$ws.Cells.Item(126, 11).Value = 3256.3638
$ws.Cells.Item(126, 12).Value = 4157.5002
$ws.Cells.Item(126, 13).Value = -786.3638000000001
$ws.Cells.Item(126, 14).Value = -9097.5002
$ws.Cells.Item(132, 8).Value = 2301.5386
$ws.Cells.Item(132, 9).Value = 1592
$ws.Cells.Item(132, 10).Value = 4107.636
$ws.Cells.Item(132, 11).Value = 4776
$ws.Cells.Item(132, 12).Value = 12322.908
$ws.Cells.Item(132, 13).Value = -2246
$ws.Cells.Item(132, 14).Value = -17382.908
$ws.Cells.Item(136, 8).Value = 25421.143
$ws.Cells.Item(136, 9).Value = 36741
$ws.Cells.Item(136, 10).Value = 2781.4285
$ws.Cells.Item(136, 11).Value = 110223
$ws.Cells.Item(136, 12).Value = 8344.2855
$ws.Cells.Item(136, 13).Value = -107673
$ws.Cells.Item(136, 14).Value = -13444.2855
